$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold textual data (prices/percentages formatted as
# strings, e.g. "53.960.04" or "  -1.50%  "). Force text format so Excel
# does not reinterpret these as numbers and mangle the exact formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.960.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "494.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0948"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.663.67"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "53.922.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.261.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "301.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.75"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.40"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.89"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.370"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "124.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.77"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "237.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.08"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.52%  "
